# Fixed DataProvider issue and updated DB verification test
#
# 1. Update the credentials used by the VALID / INVALID DataProvider sheets.
# 2. Add a new "EmployeeDetails" worksheet (used by the updated DB
#    verification test) with a small emp_id / emp_name table.

$wb = $excel.ActiveWorkbook

# --- 1. Update VALID sheet credentials -------------------------------------
$validSheet = $wb.Worksheets.Item("VALID")
$validSheet.Range("A2").Value = "dheeraj913"
$validSheet.Range("B2").Value = "@Dksharmais1908"

# --- 1b. Update INVALID sheet credentials -----------------------------------
$invalidSheet = $wb.Worksheets.Item("INVALID")
$invalidSheet.Range("A2").Value = "dheeraj913"
$invalidSheet.Range("B2").Value = "admindemo"

# --- 2. Add the EmployeeDetails sheet at the end of the workbook -----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "EmployeeDetails"

$ws.Range("A1").Value = "emp_id"
$ws.Range("B1").Value = "emp_name"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Dheeraj"
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Ajay"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Ravi"

# Match the formatting used by the other sheets (font style "1").
$validSheet.Range("A1:B2").Copy()
$ws.Range("A1:B4").PasteSpecial(-4122)
